$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 315, shifting existing rows 315-397 down to 316-398
$ws.Rows("315").Insert()

# Populate the newly inserted row 315 with the new record's values
$ws.Range("A315").Value = 11
$ws.Range("B315").Value = "Vega Monumental Concepción"
$ws.Range("C315").Value = "Bíobío"
$ws.Range("D315").Value = 44889
$ws.Range("E315").Value = 8
$ws.Range("F315").Value = 100114014
$ws.Range("G315").Value = "Betarraga"
$ws.Range("H315").Value = "Sin especificar"
$ws.Range("I315").Value = "Primera"
$ws.Range("J315").Value = 450
$ws.Range("K315").Value = 700
$ws.Range("L315").Value = 750
$ws.Range("M315").Value = 722
$ws.Range("N315").Value = "`$/paquete 5 unidades"
$ws.Range("O315").Value = "Región Metropolitana"
$ws.Range("P315").Value = 144
$ws.Range("Q315").Value = 5
$ws.Range("R315").Value = "Hortaliza"
